# Insert one new price-record row for "Zapallo italiano" (Macroferia Regional
# de Talca) just above the current row 509, pushing the existing rows 509-625
# down to 510-626 (dimension grows from A1:R625 to A1:R626).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 509:625 down by inserting a fresh row at 509.
$ws.Rows(509).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A509").Value = 5
$ws.Range("B509").Value = "Macroferia Regional de Talca"
$ws.Range("C509").Value = "Maule"
$ws.Range("D509").Value = 45173
$ws.Range("E509").Value = 7
$ws.Range("F509").Value = 100112032
$ws.Range("G509").Value = "Zapallo italiano"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Primera"
$ws.Range("J509").Value = 300
$ws.Range("K509").Value = 15000
$ws.Range("L509").Value = 15000
$ws.Range("M509").Value = 15000
$ws.Range("N509").Value = "$/caja 50 unidades"
$ws.Range("O509").Value = "Región de Arica y Parinacota"
$ws.Range("P509").Value = 300
$ws.Range("Q509").Value = 50
$ws.Range("R509").Value = "Hortaliza"

# Note: Rows(509).Insert() already carries the date-formatted style from the
# old row 509 (now row 510) into the new blank row, so D509 keeps its
# "yyyy-mm-dd" numeric format automatically - no extra style assignment
# needed.
